$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("uploadFirebase")

# 1. Update column D formulas for rows 2-131 to prepend the "+56" country code
for ($r = 2; $r -le 131; $r++) {
    $ws.Cells.Item($r, 4).Formula = '="+56"&' + "'Hoja 1'!F" + $r
}

# 2. Fill in columns A, B, C, E for rows 82-131 (these columns were previously
#    missing on those rows, mirroring the pattern used on rows 2-81)
for ($r = 82; $r -le 131; $r++) {
    $ws.Cells.Item($r, 1).Formula = "='Hoja 1'!B" + $r
    $ws.Cells.Item($r, 2).Formula = "='Hoja 1'!C" + $r
    $ws.Cells.Item($r, 3).Formula = "='Hoja 1'!D" + $r
    $ws.Cells.Item($r, 5).Formula = "='Hoja 1'!E" + $r
}

# 3. Fill in column F for rows 85-131 (rows 82-84 already had this column)
for ($r = 85; $r -le 131; $r++) {
    $ws.Cells.Item($r, 6).Formula = "='Hoja 1'!I" + $r
}

# 4. Create row 131's G cell (literal value, same as the rest of the column)
#    and copy formatting from row 130 so styles match exactly.
$ws.Cells.Item(131, 7).Value = "4Ms88Hw9i3okcTHI7AV6"
$ws.Range("G130").Copy()
$ws.Range("G131").PasteSpecial(-4122)

# 5. Copy number formats / styles from row 81 onto the newly populated
#    A:F cells of rows 82-131.
$src = $ws.Range("A81:F81")
for ($r = 82; $r -le 131; $r++) {
    $dst = $ws.Range("A" + $r + ":F" + $r)
    $src.Copy()
    $dst.PasteSpecial(-4122)
}

# 6. Update the "subscribedEventsId" value (shared across the whole G column,
#    rows 2-131) from a bare id to a JSON array containing that id.
for ($r = 2; $r -le 131; $r++) {
    $ws.Cells.Item($r, 7).Value = '["4Ms88Hw9i3okcTHI7AV6"]'
}

$excel.CutCopyMode = 0
